# Duplicate the "Ben Breadon" row (currently row 2) into new rows 4 and 5,
# as part of work on QR code tags.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MobileFoodDistro")

$values = @("Ben", "Breadon", "benbreadon@gmail.com", "1316", "Peach Tree Lane", "63069")
# Columns D (# in House) and F (Zip) hold digit strings that must stay text,
# so prefix them with an apostrophe to force text entry instead of a number.
$forceText = @($false, $false, $false, $true, $false, $true)

for ($r = 4; $r -le 5; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $v = $values[$c - 1]
        if ($forceText[$c - 1]) {
            $ws.Cells.Item($r, $c).Value = "'" + $v
        } else {
            $ws.Cells.Item($r, $c).Value = $v
        }
    }
}
